# feat: add 2022-Q1 data
#
# Inserts a new per-quarter holdings sheet "2022-Q1" (placed between the
# existing "2021-Q4" sheet and the running-totals "总计" sheet) and adds a
# corresponding summary row to "总计".
#
# To reproduce the sheetId numbering of the target workbook (the new
# "2022-Q1" sheet keeps the old sheetId that "总计" used to have, and a
# freshly created "总计" sheet picks up the next free id) we:
#   1. duplicate the existing "总计" sheet (the duplicate becomes the new
#      "总计", with a brand new sheetId, and keeps all of the original
#      rows completely untouched byte-for-byte);
#   2. rename/repurpose the original "总计" sheet object itself into
#      "2022-Q1" (it keeps its original sheetId) and fill it with the
#      quarter's fund holdings;
#   3. insert the new 2022-Q1 summary row at the top of the duplicated
#      totals sheet and rename it back to "总计".

$wb = $excel.ActiveWorkbook

# A per-quarter sheet whose header/row styling ("s=2" bold+border look) we
# can reuse as a formatting template via PasteSpecial (formats only).
$template = $wb.Worksheets.Item("2021-Q2")

$oldTotal = $wb.Worksheets.Item("总计")

# Step 1: duplicate "总计" right after itself -- the copy will carry on as
# the new "总计" sheet and keeps every original row/value untouched.
$oldTotal.Copy($null, $oldTotal)
$newTotal = $wb.Worksheets.Item($wb.Worksheets.Count)

function Set-TextCell($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)  # xlPasteFormats
}

# ---- Step 2: repurpose the original "总计" sheet into "2022-Q1" --------

$oldTotal.Cells.Clear()
$oldTotal.Name = "2022-Q1"
$q1 = $oldTotal

# Header row
Set-TextCell $q1.Range("B1") "基金代码"
Set-TextCell $q1.Range("C1") "基金名称"
Set-TextCell $q1.Range("D1") "基金规模"
Set-TextCell $q1.Range("E1") "股票总仓位"
Set-TextCell $q1.Range("F1") "仓位占比"
Set-TextCell $q1.Range("G1") "持有市值(亿元)"
Set-TextCell $q1.Range("H1") "仓位排名"
Copy-Format $template.Range("B1:H1") $q1.Range("B1:H1")

# Data rows
$q1.Range("A2").Value = 0
Set-TextCell $q1.Range("B2") "006923"
Set-TextCell $q1.Range("C2") "前海开源沪港深非周期性行业股票A"
Set-TextCell $q1.Range("D2") "0.54"
Set-TextCell $q1.Range("E2") "93.77"
Set-TextCell $q1.Range("F2") "4.81"
Set-TextCell $q1.Range("G2") "0.0260"
$q1.Range("H2").Value = 10

$q1.Range("A3").Value = 1
Set-TextCell $q1.Range("B3") "006924"
Set-TextCell $q1.Range("C3") "前海开源沪港深非周期性行业股票C"
Set-TextCell $q1.Range("D3") "0.22"
Set-TextCell $q1.Range("E3") "93.77"
Set-TextCell $q1.Range("F3") "4.81"
Set-TextCell $q1.Range("G3") "0.0106"
$q1.Range("H3").Value = 10

$q1.Range("A4").Value = 2
Set-TextCell $q1.Range("B4") "004098"
Set-TextCell $q1.Range("C4") "前海开源港股通股息率50强股票"
Set-TextCell $q1.Range("D4") "0.34"
Set-TextCell $q1.Range("E4") "88.92"
Set-TextCell $q1.Range("F4") "2.41"
Set-TextCell $q1.Range("G4") "0.0082"
$q1.Range("H4").Value = 8

Copy-Format $template.Range("A2") $q1.Range("A2:A4")

# ---- Step 3: insert the new row into the duplicated totals sheet -------

$newTotal.Range("A2:D2").Insert(-4121)  # xlShiftDown

# Insert() makes the new row inherit some formatting from the row below;
# start from a clean slate like the other (non-header, non-index) data
# cells in this workbook.
$newTotal.Range("A2:D2").ClearFormats()

$newTotal.Range("A2").Value = 0
Set-TextCell $newTotal.Range("B2") "2022-Q1"
$newTotal.Range("C2").Value = 3
$newTotal.Range("D2").Value = 0.04
Copy-Format $template.Range("A2") $newTotal.Range("A2")

# The row-0/1/2/3 index column needs renumbering since the previously
# existing rows shifted down by one.
$newTotal.Range("A3").Value = 1
$newTotal.Range("A4").Value = 2
$newTotal.Range("A5").Value = 3
$newTotal.Range("A6").Value = 4

$newTotal.Name = "总计"
